$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on cells whose new values would otherwise be
# auto-converted to numbers by Excel, so they stay text like the source file.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range('D2').Value = '41.950.51'
$ws.Range('E2').Value = '  +5.02%  '
$ws.Range('D3').Value = '2.255.20'
$ws.Range('E3').Value = '  +1.57%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '301.96'
$ws.Range('E5').Value = '  +3.51%  '
$ws.Range('D6').Value = '92.79'
$ws.Range('E6').Value = '  +5.94%  '
$ws.Range('D7').Value = '0.532'
$ws.Range('E7').Value = '  +3.65%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').Value = '0.483'
$ws.Range('E9').Value = '  +3.38%  '
$ws.Range('D10').Value = '32.71'
$ws.Range('E10').Value = '  +7.23%  '
$ws.Range('D11').Value = '54.59'
$ws.Range('E11').Value = '  +9.20%  '
$ws.Range('D12').Value = '0.0801'
$ws.Range('E12').Value = '  +2.53%  '
$ws.Range('E13').Value = '  +3.29%  '
$ws.Range('D14').Value = '6.69'
$ws.Range('E14').Value = '  +3.89%  '
$ws.Range('D15').Value = '2.604.28'
$ws.Range('E15').Value = '  +1.65%  '
$ws.Range('D16').Value = '14.14'
$ws.Range('E16').Value = '  +2.48%  '
$ws.Range('D17').Value = '2.271.61'
$ws.Range('E17').Value = '  +2.80%  '
$ws.Range('D18').Value = '0.757'
$ws.Range('E18').Value = '  +3.54%  '
$ws.Range('D19').Value = '41.856.38'
$ws.Range('E19').Value = '  +4.96%  '
$ws.Range('D20').Value = '12.27'
$ws.Range('E20').Value = '  +10.49%  '
$ws.Range('D21').Value = '0.0₃0904'
$ws.Range('E21').Value = '  +1.99%  '
$ws.Range('D22').Value = '5.95'
$ws.Range('E22').Value = '  +3.67%  '
$ws.Range('D23').Value = '67.20'
$ws.Range('E23').Value = '  +2.31%  '
$ws.Range('D24').Value = '241.92'
$ws.Range('E24').Value = '  +1.88%  '
$ws.Range('E25').Value = '  +5.40%  '
$ws.Range('D26').Value = '1.93'
$ws.Range('E26').Value = '  +5.14%  '
$ws.Range('E27').Value = '  -0.13%  '
$ws.Range('D28').Value = '23.96'
$ws.Range('E28').Value = '  +3.57%  '
$ws.Range('E29').Value = '  +6.22%  '
$ws.Range('D30').Value = '9.68'
$ws.Range('E30').Value = '  +4.87%  '
$ws.Range('D31').Value = '34.15'
$ws.Range('E31').Value = '  +6.62%  '
$ws.Range('D32').Value = '158.72'
$ws.Range('E32').Value = '  +1.40%  '
$ws.Range('D33').Value = '1.00'
$ws.Range('E33').Value = '  +0.06%  '
$ws.Range('D34').Value = '5.16'
$ws.Range('E34').Value = '  +3.96%  '
$ws.Range('D35').Value = '0.0744'
$ws.Range('E35').Value = '  +4.54%  '
$ws.Range('D36').Value = '3.04'
$ws.Range('E36').Value = '  +2.13%  '
$ws.Range('E37').Value = '  +2.74%  '
$ws.Range('B38').Value = 'Kaspa'
$ws.Range('C38').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D38').Value = '0.105'
$ws.Range('E38').Value = '  +5.84%  '
$ws.Range('B39').Value = 'Celestia'
$ws.Range('C39').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D39').Value = '16.67'
$ws.Range('E39').Value = '  +8.67%  '
$ws.Range('E40').Value = '  +3.74%  '
$ws.Range('D41').Value = '1.79'
$ws.Range('E41').Value = '  +3.53%  '
$ws.Range('D42').Value = '3.93'
$ws.Range('E42').Value = '  +5.73%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '2.050.54'
$ws.Range('E43').Value = '  -2.76%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').Value = '19.82'
$ws.Range('E44').Value = '  +9.39%  '
$ws.Range('D45').Value = '0.0279'
$ws.Range('E45').Value = '  +3.38%  '
$ws.Range('D46').Value = '10.13'
$ws.Range('E46').Value = '  +2.48%  '
$ws.Range('D47').Value = '2.87'
$ws.Range('E47').Value = '  +5.92%  '
$ws.Range('D48').Value = '2.05'
$ws.Range('E48').Value = '  +1.89%  '
$ws.Range('D49').Value = '1.52'
$ws.Range('E49').Value = '  +3.49%  '
$ws.Range('D50').Value = '1.14'
$ws.Range('E50').Value = '  +3.90%  '
$ws.Range('D51').Value = '51.82'
$ws.Range('E51').Value = '  +5.73%  '
